$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff (189 changed cells across rows 2-28).
# Each odds/value cell is updated to its new numeric value; SnapshotTS (BH) cells get the new timestamp string.

# Row 2
$ws.Range("J2").Value = 2.98
$ws.Range("BH2").Value = "2026-02-23 12:51:28"

# Row 3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.45
$ws.Range("J3").Value = 4
$ws.Range("N3").Value = 4.7
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2.28
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.86
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 2.42
$ws.Range("X3").Value = 20
$ws.Range("Y3").Value = 16
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AE3").Value = 40
$ws.Range("AF3").Value = 15.5
$ws.Range("AI3").Value = 48
$ws.Range("AJ3").Value = 27
$ws.Range("AN3").Value = 13
$ws.Range("AP3").Value = 17.5
$ws.Range("AS3").Value = 44
$ws.Range("AU3").Value = 8.4
$ws.Range("AW3").Value = 30
$ws.Range("BA3").Value = 34
$ws.Range("BB3").Value = 23
$ws.Range("BC3").Value = 18.5
$ws.Range("BE3").Value = 32
$ws.Range("BG3").Value = 24
$ws.Range("BH3").Value = "2026-02-23 12:51:28"

# Row 4
$ws.Range("F4").Value = 2.72
$ws.Range("G4").Value = 2.92
$ws.Range("I4").Value = 3
$ws.Range("K4").Value = 3.45
$ws.Range("BH4").Value = "2026-02-23 12:51:28"

# Row 5
$ws.Range("BH5").Value = "2026-02-23 12:51:28"

# Row 6
$ws.Range("BH6").Value = "2026-02-23 12:51:28"

# Row 7
$ws.Range("BH7").Value = "2026-02-23 12:51:28"

# Row 8
$ws.Range("BH8").Value = "2026-02-23 12:51:28"

# Row 9
$ws.Range("K9").Value = 10.5
$ws.Range("P9").Value = 4.1
$ws.Range("Q9").Value = 1.22
$ws.Range("BH9").Value = "2026-02-23 12:51:28"

# Row 10
$ws.Range("G10").Value = 3.1
$ws.Range("Q10").Value = 2.08
$ws.Range("BH10").Value = "2026-02-23 12:51:28"

# Row 11
$ws.Range("P11").Value = 1.43
$ws.Range("Q11").Value = 2.92
$ws.Range("BH11").Value = "2026-02-23 12:51:28"

# Row 12
$ws.Range("G12").Value = 1.76
$ws.Range("H12").Value = 4.9
$ws.Range("Q12").Value = 2.58
$ws.Range("BH12").Value = "2026-02-23 12:51:28"

# Row 13
$ws.Range("G13").Value = 2.82
$ws.Range("I13").Value = 4.1
$ws.Range("J13").Value = 2.76
$ws.Range("P13").Value = 1.62
$ws.Range("Q13").Value = 2.2
$ws.Range("T13").Value = 1.94
$ws.Range("U13").Value = 1.84
$ws.Range("X13").Value = 12
$ws.Range("AC13").Value = 8.800000000000001
$ws.Range("AE13").Value = 65
$ws.Range("AJ13").Value = 44
$ws.Range("AK13").Value = 40
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 180
$ws.Range("AN13").Value = 38
$ws.Range("AP13").Value = 3.2
$ws.Range("AQ13").Value = 3.3
$ws.Range("AR13").Value = 3.8
$ws.Range("AS13").Value = 4.2
$ws.Range("AT13").Value = 7
$ws.Range("AV13").Value = 3.55
$ws.Range("AW13").Value = 4.1
$ws.Range("AX13").Value = 3.5
$ws.Range("AY13").Value = 3.35
$ws.Range("AZ13").Value = 3.7
$ws.Range("BA13").Value = 4.1
$ws.Range("BB13").Value = 4
$ws.Range("BC13").Value = 3.95
$ws.Range("BD13").Value = 4.1
$ws.Range("BE13").Value = 4.3
$ws.Range("BF13").Value = 3.95
$ws.Range("BG13").Value = 4.1
$ws.Range("BH13").Value = "2026-02-23 12:51:28"

# Row 14
$ws.Range("F14").Value = 2.54
$ws.Range("G14").Value = 2.92
$ws.Range("I14").Value = 2.96
$ws.Range("J14").Value = 3.75
$ws.Range("BH14").Value = "2026-02-23 12:51:28"

# Row 15
$ws.Range("F15").Value = 3.55
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 1.89
$ws.Range("I15").Value = 2.06
$ws.Range("K15").Value = 4.7
$ws.Range("BH15").Value = "2026-02-23 12:51:28"

# Row 16
$ws.Range("H16").Value = 2.84
$ws.Range("I16").Value = 3.45
$ws.Range("J16").Value = 3.25
$ws.Range("Q16").Value = 1.92
$ws.Range("BH16").Value = "2026-02-23 12:51:28"

# Row 17
$ws.Range("G17").Value = 2.6
$ws.Range("K17").Value = 3.95
$ws.Range("P17").Value = 1.96
$ws.Range("Q17").Value = 1.84
$ws.Range("BH17").Value = "2026-02-23 12:51:28"

# Row 18
$ws.Range("H18").Value = 3.05
$ws.Range("K18").Value = 4.3
$ws.Range("P18").Value = 1.8
$ws.Range("Q18").Value = 1.98
$ws.Range("BH18").Value = "2026-02-23 12:51:28"

# Row 19
$ws.Range("P19").Value = 1.36
$ws.Range("BH19").Value = "2026-02-23 12:51:28"

# Row 20
$ws.Range("H20").Value = 1.44
$ws.Range("BH20").Value = "2026-02-23 12:51:28"

# Row 21
$ws.Range("F21").Value = 1.53
$ws.Range("H21").Value = 6
$ws.Range("J21").Value = 5.4
$ws.Range("K21").Value = 5.5
$ws.Range("P21").Value = 2.96
$ws.Range("S21").Value = 2.22
$ws.Range("U21").Value = 2.44
$ws.Range("X21").Value = 34
$ws.Range("Y21").Value = 32
$ws.Range("Z21").Value = 60
$ws.Range("AC21").Value = 13
$ws.Range("AD21").Value = 24
$ws.Range("AE21").Value = 85
$ws.Range("AI21").Value = 60
$ws.Range("AJ21").Value = 15
$ws.Range("AK21").Value = 14
$ws.Range("AN21").Value = 5.1
$ws.Range("AP21").Value = 27
$ws.Range("AR21").Value = 46
$ws.Range("AT21").Value = 12.5
$ws.Range("AV21").Value = 20
$ws.Range("AW21").Value = 50
$ws.Range("AZ21").Value = 17
$ws.Range("BA21").Value = 44
$ws.Range("BB21").Value = 13.5
$ws.Range("BE21").Value = 55
$ws.Range("BF21").Value = 4.9
$ws.Range("BG21").Value = 46
$ws.Range("BH21").Value = "2026-02-23 12:51:28"

# Row 22
$ws.Range("G22").Value = 1.5
$ws.Range("P22").Value = 2.88
$ws.Range("S22").Value = 2.24
$ws.Range("U22").Value = 2.38
$ws.Range("X22").Value = 30
$ws.Range("AG22").Value = 10.5
$ws.Range("AM22").Value = 75
$ws.Range("AN22").Value = 4.9
$ws.Range("AP22").Value = 27
$ws.Range("AQ22").Value = 30
$ws.Range("AR22").Value = 55
$ws.Range("AU22").Value = 11.5
$ws.Range("AV22").Value = 24
$ws.Range("AW22").Value = 60
$ws.Range("AZ22").Value = 17.5
$ws.Range("BA22").Value = 55
$ws.Range("BD22").Value = 23
$ws.Range("BE22").Value = 60
$ws.Range("BF22").Value = 4.7
$ws.Range("BG22").Value = 55
$ws.Range("BH22").Value = "2026-02-23 12:51:28"

# Row 23
$ws.Range("T23").Value = 1.69
$ws.Range("U23").Value = 2.36
$ws.Range("Z23").Value = 120
$ws.Range("AA23").Value = 390
$ws.Range("AC23").Value = 18
$ws.Range("AG23").Value = 11.5
$ws.Range("AS23").Value = 110
$ws.Range("AT23").Value = 15
$ws.Range("BG23").Value = 38
$ws.Range("BH23").Value = "2026-02-23 12:51:28"

# Row 24
$ws.Range("BH24").Value = "2026-02-23 12:51:28"

# Row 25
$ws.Range("BH25").Value = "2026-02-23 12:51:28"

# Row 26
$ws.Range("F26").Value = 1.89
$ws.Range("I26").Value = 7.4
$ws.Range("P26").Value = 1.43
$ws.Range("Q26").Value = 2.9
$ws.Range("BH26").Value = "2026-02-23 12:51:28"

# Row 27
$ws.Range("G27").Value = 3.85
$ws.Range("J27").Value = 2.78
$ws.Range("Q27").Value = 2.98
$ws.Range("BH27").Value = "2026-02-23 12:51:28"

# Row 28
$ws.Range("BH28").Value = "2026-02-23 12:51:28"
